# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same shape as the
# other quarterly sheets) right before the "总计" (totals) sheet, and adds a
# corresponding summary row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Borrow formatting from an existing quarterly sheet (same columns/style)
$fmtSource = $wb.Worksheets.Item("2021-Q4")
$fmtSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSource.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 003293 易方达科瑞灵活配置混合
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'003293"
$q1.Range("C2").Value = "易方达科瑞灵活配置混合"
$q1.Range("D2").Value = "'34.67"
$q1.Range("E2").Value = "'78.17"
$q1.Range("F2").Value = "'2.31"
$q1.Range("G2").Value = "'0.8009"
$q1.Range("H2").Value = 9

# Row 3 - 110012 易方达科汇灵活配置混合
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'110012"
$q1.Range("C3").Value = "易方达科汇灵活配置混合"
$q1.Range("D3").Value = "'15.73"
$q1.Range("E3").Value = "'75.64"
$q1.Range("F3").Value = "'2.50"
$q1.Range("G3").Value = "'0.3932"
$q1.Range("H3").Value = 9

# Row 4 - 011649 易方达逆向投资混合A
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'011649"
$q1.Range("C4").Value = "易方达逆向投资混合A"
$q1.Range("D4").Value = "'7.49"
$q1.Range("E4").Value = "'85.02"
$q1.Range("F4").Value = "'3.11"
$q1.Range("G4").Value = "'0.2329"
$q1.Range("H4").Value = 8

# Row 5 - 011650 易方达逆向投资混合C
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'011650"
$q1.Range("C5").Value = "易方达逆向投资混合C"
$q1.Range("D5").Value = "'1.96"
$q1.Range("E5").Value = "'85.02"
$q1.Range("F5").Value = "'3.11"
$q1.Range("G5").Value = "'0.0610"
$q1.Range("H5").Value = 8

# Row 6 - 540004 汇丰晋信2026周期混合
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "'540004"
$q1.Range("C6").Value = "汇丰晋信2026周期混合"
$q1.Range("D6").Value = "'1.14"
$q1.Range("E6").Value = "'31.29"
$q1.Range("F6").Value = "'1.75"
$q1.Range("G6").Value = "'0.0200"
$q1.Range("H6").Value = 6

# ---------------------------------------------------------------------
# 2) Prepend a 2022-Q1 summary row into "总计"
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("A2:D2").ClearFormats()

# copy the number-column (A) style down from the row it pushed out
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 5
$tot.Range("D2").Value = 1.51

# renumber the existing index column, which Insert() does not do for us
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
